$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Select the source range on Sheet1 (becomes Sheet1's remembered selection
# state) and copy it so it can be pasted into Sheet2 further right.
$ws1.Activate()
$ws1.Range("A1:A7").Select()
$ws1.Range("A1:A7").Copy()

# The single-column list on Sheet2 was wrongly anchored at column A;
# insert 3 blank columns before it so it shifts from A:A into D:D.
$ws2.Activate()
$ws2.Range("A:C").Insert()

# Paste the Sheet1 range into Sheet2 next to the (now-shifted) list.
$ws2.Range("N4").PasteSpecial()
$excel.CutCopyMode = $false

# Leave Sheet2 active with the fixed single-column range selected.
$ws2.Range("D2:D17").Select()
